$p = $ppt.ActivePresentation

# 1. Table on slide 16 switches from the custom "Table_0" style to the
#    built-in "Medium Style 2 - Accent 1" table style.
$s = $p.Slides.Item(16)
$tbl = $s.Shapes.Item(3).Table
$tbl.ApplyStyle("{C53F6008-7C91-41F1-B57A-9E2306D5BF75}")

# 2. The deck's main theme (used by the slide master / all slides) is
#    switched from the "Integral" color scheme to the standard
#    "Office" color scheme.
$sr = $p.Slides.Range()
$tcs = $sr.ThemeColorScheme
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
